$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7: the "ser:" numbers shift up - blog post 68 retires, 69 and 71 go live.
$ws.Range("G7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 69"
$ws.Range("C7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 71"
